$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 17
$ws.Range("B2").Value = 'memory'
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 161
$ws.Range("G2").Value = 'living_rooms'
$ws.Range("H2").Value = 'living_rooms'
$ws.Range("I2").Value = 'target'
$ws.Range("J2").Value = 'old'
$ws.Range("K2").Value = 'j'
$ws.Range("L2").Value = 'stimuli/img_6zz63.png'
$ws.Range("M2").Value = 87.66666666666667
$ws.Range("N2").Value = 70.6
$ws.Range("O2").Value = 79.13333333333333
$ws.Range("P2").Value = 45
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 10
$ws.Range("A3").Value = 17
$ws.Range("B3").Value = 'memory'
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 162
$ws.Range("G3").Value = 'living_rooms'
$ws.Range("H3").Value = 'living_rooms'
$ws.Range("I3").Value = $null
$ws.Range("J3").Value = 'new'
$ws.Range("K3").Value = 'f'
$ws.Range("L3").Value = 'stimuli/img_b17ma.png'
$ws.Range("M3").Value = 23.0625
$ws.Range("N3").Value = 13.375
$ws.Range("O3").Value = 18.21875
$ws.Range("P3").Value = 48
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("A4").Value = 17
$ws.Range("B4").Value = 'memory'
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 163
$ws.Range("G4").Value = 'living_rooms'
$ws.Range("H4").Value = 'living_rooms'
$ws.Range("I4").Value = 'target'
$ws.Range("J4").Value = 'old'
$ws.Range("K4").Value = 'j'
$ws.Range("L4").Value = 'stimuli/img_wgkqa.png'
$ws.Range("M4").Value = 87.25581395348837
$ws.Range("N4").Value = 71.13953488372093
$ws.Range("O4").Value = 79.19767441860465
$ws.Range("P4").Value = 43
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = 10
$ws.Range("S4").Value = 10
$ws.Range("A5").Value = 17
$ws.Range("B5").Value = 'memory'
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 164
$ws.Range("G5").Value = 'living_rooms'
$ws.Range("H5").Value = 'living_rooms'
$ws.Range("I5").Value = $null
$ws.Range("J5").Value = 'new'
$ws.Range("K5").Value = 'f'
$ws.Range("L5").Value = 'stimuli/img_3jnt7.png'
$ws.Range("M5").Value = 49.52272727272727
$ws.Range("N5").Value = 35.25
$ws.Range("O5").Value = 42.38636363636364
$ws.Range("P5").Value = 44
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 3
$ws.Range("A6").Value = 17
$ws.Range("B6").Value = 'memory'
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 165
$ws.Range("G6").Value = 'living_rooms'
$ws.Range("H6").Value = 'living_rooms'
$ws.Range("I6").Value = 'target'
$ws.Range("J6").Value = 'old'
$ws.Range("K6").Value = 'j'
$ws.Range("L6").Value = 'stimuli/img_amsgw.png'
$ws.Range("M6").Value = 86.08510638297872
$ws.Range("N6").Value = 65.95744680851064
$ws.Range("O6").Value = 76.02127659574468
$ws.Range("P6").Value = 47
$ws.Range("Q6").Value = 9
$ws.Range("R6").Value = 9
$ws.Range("S6").Value = 9
$ws.Range("A7").Value = 17
$ws.Range("B7").Value = 'memory'
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 166
$ws.Range("G7").Value = 'living_rooms'
$ws.Range("H7").Value = 'living_rooms'
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = 'new'
$ws.Range("K7").Value = 'f'
$ws.Range("L7").Value = 'stimuli/img_15bss.png'
$ws.Range("M7").Value = 88.42222222222222
$ws.Range("N7").Value = 75.35555555555555
$ws.Range("O7").Value = 81.88888888888889
$ws.Range("P7").Value = 45
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = 10
$ws.Range("S7").Value = 10
$ws.Range("A8").Value = 17
$ws.Range("B8").Value = 'memory'
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 167
$ws.Range("G8").Value = 'living_rooms'
$ws.Range("H8").Value = 'living_rooms'
$ws.Range("I8").Value = 'target'
$ws.Range("J8").Value = 'old'
$ws.Range("K8").Value = 'j'
$ws.Range("L8").Value = 'stimuli/img_eh0no.png'
$ws.Range("M8").Value = 53.66666666666666
$ws.Range("N8").Value = 36.02564102564103
$ws.Range("O8").Value = 44.84615384615385
$ws.Range("P8").Value = 39
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 3
$ws.Range("S8").Value = 3
$ws.Range("A9").Value = 17
$ws.Range("B9").Value = 'memory'
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 168
$ws.Range("G9").Value = 'living_rooms'
$ws.Range("H9").Value = 'living_rooms'
$ws.Range("I9").Value = $null
$ws.Range("J9").Value = 'new'
$ws.Range("K9").Value = 'f'
$ws.Range("L9").Value = 'stimuli/img_314bq.png'
$ws.Range("M9").Value = 37.08888888888889
$ws.Range("N9").Value = 20.04444444444444
$ws.Range("O9").Value = 28.56666666666667
$ws.Range("P9").Value = 45
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2
$ws.Range("A10").Value = 17
$ws.Range("B10").Value = 'memory'
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 169
$ws.Range("G10").Value = 'living_rooms'
$ws.Range("H10").Value = 'living_rooms'
$ws.Range("I10").Value = 'target'
$ws.Range("J10").Value = 'old'
$ws.Range("K10").Value = 'j'
$ws.Range("L10").Value = 'stimuli/img_xbtev.png'
$ws.Range("M10").Value = 13.68181818181818
$ws.Range("N10").Value = 8.568181818181818
$ws.Range("O10").Value = 11.125
$ws.Range("P10").Value = 44
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("A11").Value = 17
$ws.Range("B11").Value = 'memory'
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 170
$ws.Range("G11").Value = 'living_rooms'
$ws.Range("H11").Value = 'living_rooms'
$ws.Range("I11").Value = $null
$ws.Range("J11").Value = 'new'
$ws.Range("K11").Value = 'f'
$ws.Range("L11").Value = 'stimuli/img_dg5h7.png'
$ws.Range("M11").Value = 88.72093023255815
$ws.Range("N11").Value = 76.06976744186046
$ws.Range("O11").Value = 82.3953488372093
$ws.Range("P11").Value = 43
$ws.Range("Q11").Value = 10
$ws.Range("R11").Value = 10
$ws.Range("S11").Value = 10
$ws.Range("A12").Value = 17
$ws.Range("B12").Value = 'memory'
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 171
$ws.Range("G12").Value = 'living_rooms'
$ws.Range("H12").Value = 'living_rooms'
$ws.Range("I12").Value = 'target'
$ws.Range("J12").Value = 'old'
$ws.Range("K12").Value = 'j'
$ws.Range("L12").Value = 'stimuli/img_xu1p3.png'
$ws.Range("M12").Value = 75.27659574468085
$ws.Range("N12").Value = 56.68085106382978
$ws.Range("O12").Value = 65.97872340425532
$ws.Range("P12").Value = 47
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("A13").Value = 17
$ws.Range("B13").Value = 'memory'
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = 172
$ws.Range("G13").Value = 'living_rooms'
$ws.Range("H13").Value = 'living_rooms'
$ws.Range("I13").Value = 'target'
$ws.Range("J13").Value = 'old'
$ws.Range("K13").Value = 'j'
$ws.Range("L13").Value = 'stimuli/img_cehin.png'
$ws.Range("M13").Value = 78.86363636363636
$ws.Range("N13").Value = 60.02272727272727
$ws.Range("O13").Value = 69.44318181818181
$ws.Range("P13").Value = 44
$ws.Range("Q13").Value = 7
$ws.Range("R13").Value = 7
$ws.Range("S13").Value = 7
$ws.Range("A14").Value = 17
$ws.Range("B14").Value = 'memory'
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 173
$ws.Range("G14").Value = 'living_rooms'
$ws.Range("H14").Value = 'living_rooms'
$ws.Range("I14").Value = $null
$ws.Range("J14").Value = 'new'
$ws.Range("K14").Value = 'f'
$ws.Range("L14").Value = 'stimuli/img_0jzz7.png'
$ws.Range("M14").Value = 84.85106382978724
$ws.Range("N14").Value = 68.87234042553192
$ws.Range("O14").Value = 76.86170212765958
$ws.Range("P14").Value = 47
$ws.Range("Q14").Value = 9
$ws.Range("R14").Value = 9
$ws.Range("S14").Value = 9
$ws.Range("A15").Value = 17
$ws.Range("B15").Value = 'memory'
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = 174
$ws.Range("G15").Value = 'living_rooms'
$ws.Range("H15").Value = 'living_rooms'
$ws.Range("I15").Value = 'target'
$ws.Range("J15").Value = 'old'
$ws.Range("K15").Value = 'j'
$ws.Range("L15").Value = 'stimuli/img_w8yhd.png'
$ws.Range("M15").Value = 55.74418604651163
$ws.Range("N15").Value = 38.90697674418605
$ws.Range("O15").Value = 47.32558139534883
$ws.Range("P15").Value = 43
$ws.Range("Q15").Value = 4
$ws.Range("R15").Value = 4
$ws.Range("S15").Value = 4
$ws.Range("A16").Value = 17
$ws.Range("B16").Value = 'memory'
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 175
$ws.Range("G16").Value = 'living_rooms'
$ws.Range("H16").Value = 'living_rooms'
$ws.Range("I16").Value = 'target'
$ws.Range("J16").Value = 'old'
$ws.Range("K16").Value = 'j'
$ws.Range("L16").Value = 'stimuli/img_16kib.png'
$ws.Range("M16").Value = 80.97727272727273
$ws.Range("N16").Value = 61.11363636363637
$ws.Range("O16").Value = 71.04545454545455
$ws.Range("P16").Value = 44
$ws.Range("Q16").Value = 8
$ws.Range("R16").Value = 8
$ws.Range("S16").Value = 8
$ws.Range("A17").Value = 17
$ws.Range("B17").Value = 'memory'
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 176
$ws.Range("G17").Value = 'living_rooms'
$ws.Range("H17").Value = 'living_rooms'
$ws.Range("I17").Value = $null
$ws.Range("J17").Value = 'new'
$ws.Range("K17").Value = 'f'
$ws.Range("L17").Value = 'stimuli/img_5mw7y.png'
$ws.Range("M17").Value = 72.6590909090909
$ws.Range("N17").Value = 50.86363636363637
$ws.Range("O17").Value = 61.76136363636364
$ws.Range("P17").Value = 44
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = 6
$ws.Range("S17").Value = 6
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 'memory'
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 17
$ws.Range("F18").Value = 177
$ws.Range("G18").Value = 'living_rooms'
$ws.Range("H18").Value = 'living_rooms'
$ws.Range("I18").Value = 'target'
$ws.Range("J18").Value = 'old'
$ws.Range("K18").Value = 'j'
$ws.Range("L18").Value = 'stimuli/img_di6f0.png'
$ws.Range("M18").Value = 94.04347826086956
$ws.Range("N18").Value = 83.34782608695652
$ws.Range("O18").Value = 88.69565217391303
$ws.Range("P18").Value = 46
$ws.Range("Q18").Value = 10
$ws.Range("R18").Value = 10
$ws.Range("S18").Value = 10
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 'memory'
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 18
$ws.Range("F19").Value = 178
$ws.Range("G19").Value = 'living_rooms'
$ws.Range("H19").Value = 'living_rooms'
$ws.Range("I19").Value = 'target'
$ws.Range("J19").Value = 'old'
$ws.Range("K19").Value = 'j'
$ws.Range("L19").Value = 'stimuli/img_kost0.png'
$ws.Range("M19").Value = 63.09090909090909
$ws.Range("N19").Value = 42.77272727272727
$ws.Range("O19").Value = 52.93181818181819
$ws.Range("P19").Value = 44
$ws.Range("Q19").Value = 5
$ws.Range("R19").Value = 5
$ws.Range("S19").Value = 5
$ws.Range("A20").Value = 17
$ws.Range("B20").Value = 'memory'
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 19
$ws.Range("F20").Value = 179
$ws.Range("G20").Value = 'living_rooms'
$ws.Range("H20").Value = $null
$ws.Range("I20").Value = $null
$ws.Range("J20").Value = 'catch'
$ws.Range("K20").Value = 'f'
$ws.Range("L20").Value = 'stimuli/catch_02.jpg'
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = $null
$ws.Range("O20").Value = $null
$ws.Range("P20").Value = $null
$ws.Range("Q20").Value = $null
$ws.Range("R20").Value = $null
$ws.Range("S20").Value = $null
$ws.Range("A21").Value = 17
$ws.Range("B21").Value = 'memory'
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 180
$ws.Range("G21").Value = 'living_rooms'
$ws.Range("H21").Value = 'living_rooms'
$ws.Range("I21").Value = $null
$ws.Range("J21").Value = 'new'
$ws.Range("K21").Value = 'f'
$ws.Range("L21").Value = 'stimuli/img_zh8ms.png'
$ws.Range("M21").Value = 59.82608695652174
$ws.Range("N21").Value = 39.43478260869565
$ws.Range("O21").Value = 49.6304347826087
$ws.Range("P21").Value = 46
$ws.Range("Q21").Value = 4
$ws.Range("R21").Value = 4
$ws.Range("S21").Value = 4
$ws.Range("A22").Value = 17
$ws.Range("B22").Value = 'memory'
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 21
$ws.Range("F22").Value = 181
$ws.Range("G22").Value = 'living_rooms'
$ws.Range("H22").Value = 'living_rooms'
$ws.Range("I22").Value = $null
$ws.Range("J22").Value = 'new'
$ws.Range("K22").Value = 'f'
$ws.Range("L22").Value = 'stimuli/img_f63yi.png'
$ws.Range("M22").Value = 85.275
$ws.Range("N22").Value = 68.475
$ws.Range("O22").Value = 76.875
$ws.Range("P22").Value = 40
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9
$ws.Range("A23").Value = 17
$ws.Range("B23").Value = 'memory'
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 22
$ws.Range("F23").Value = 182
$ws.Range("G23").Value = 'living_rooms'
$ws.Range("H23").Value = 'living_rooms'
$ws.Range("I23").Value = $null
$ws.Range("J23").Value = 'new'
$ws.Range("K23").Value = 'f'
$ws.Range("L23").Value = 'stimuli/img_eiu3c.png'
$ws.Range("M23").Value = 65.1590909090909
$ws.Range("N23").Value = 46.22727272727273
$ws.Range("O23").Value = 55.69318181818181
$ws.Range("P23").Value = 44
$ws.Range("Q23").Value = 5
$ws.Range("R23").Value = 5
$ws.Range("S23").Value = 5
$ws.Range("A24").Value = 17
$ws.Range("B24").Value = 'memory'
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 23
$ws.Range("F24").Value = 183
$ws.Range("G24").Value = 'living_rooms'
$ws.Range("H24").Value = 'living_rooms'
$ws.Range("I24").Value = $null
$ws.Range("J24").Value = 'new'
$ws.Range("K24").Value = 'f'
$ws.Range("L24").Value = 'stimuli/img_s2zoe.png'
$ws.Range("M24").Value = 64.71428571428571
$ws.Range("N24").Value = 44.90476190476191
$ws.Range("O24").Value = 54.80952380952381
$ws.Range("P24").Value = 42
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5
$ws.Range("A25").Value = 17
$ws.Range("B25").Value = 'memory'
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 24
$ws.Range("F25").Value = 184
$ws.Range("G25").Value = 'living_rooms'
$ws.Range("H25").Value = 'living_rooms'
$ws.Range("I25").Value = $null
$ws.Range("J25").Value = 'new'
$ws.Range("K25").Value = 'f'
$ws.Range("L25").Value = 'stimuli/img_qrc78.png'
$ws.Range("M25").Value = 76.2
$ws.Range("N25").Value = 59.875
$ws.Range("O25").Value = 68.0375
$ws.Range("P25").Value = 40
$ws.Range("Q25").Value = 7
$ws.Range("R25").Value = 7
$ws.Range("S25").Value = 7
$ws.Range("A26").Value = 17
$ws.Range("B26").Value = 'memory'
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 185
$ws.Range("G26").Value = 'living_rooms'
$ws.Range("H26").Value = 'living_rooms'
$ws.Range("I26").Value = 'target'
$ws.Range("J26").Value = 'old'
$ws.Range("K26").Value = 'j'
$ws.Range("L26").Value = 'stimuli/img_bj99b.png'
$ws.Range("M26").Value = 82.79069767441861
$ws.Range("N26").Value = 65.46511627906976
$ws.Range("O26").Value = 74.12790697674419
$ws.Range("P26").Value = 43
$ws.Range("Q26").Value = 8
$ws.Range("R26").Value = 8
$ws.Range("S26").Value = 8
$ws.Range("A27").Value = 17
$ws.Range("B27").Value = 'memory'
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 26
$ws.Range("F27").Value = 186
$ws.Range("G27").Value = 'living_rooms'
$ws.Range("H27").Value = 'living_rooms'
$ws.Range("I27").Value = 'target'
$ws.Range("J27").Value = 'old'
$ws.Range("K27").Value = 'j'
$ws.Range("L27").Value = 'stimuli/img_pey7u.png'
$ws.Range("M27").Value = 30.34883720930232
$ws.Range("N27").Value = 20.34883720930232
$ws.Range("O27").Value = 25.34883720930232
$ws.Range("P27").Value = 43
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 2
$ws.Range("S27").Value = 2
$ws.Range("A28").Value = 17
$ws.Range("B28").Value = 'memory'
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 27
$ws.Range("F28").Value = 187
$ws.Range("G28").Value = 'living_rooms'
$ws.Range("H28").Value = 'living_rooms'
$ws.Range("I28").Value = 'target'
$ws.Range("J28").Value = 'old'
$ws.Range("K28").Value = 'j'
$ws.Range("L28").Value = 'stimuli/img_bbs77.png'
$ws.Range("M28").Value = 31.64444444444445
$ws.Range("N28").Value = 21.26666666666667
$ws.Range("O28").Value = 26.45555555555556
$ws.Range("P28").Value = 45
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 2
$ws.Range("S28").Value = 2
$ws.Range("A29").Value = 17
$ws.Range("B29").Value = 'memory'
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 28
$ws.Range("F29").Value = 188
$ws.Range("G29").Value = 'living_rooms'
$ws.Range("H29").Value = 'living_rooms'
$ws.Range("I29").Value = $null
$ws.Range("J29").Value = 'new'
$ws.Range("K29").Value = 'f'
$ws.Range("L29").Value = 'stimuli/img_il020.png'
$ws.Range("M29").Value = 18.85416666666667
$ws.Range("N29").Value = 16.16666666666667
$ws.Range("O29").Value = 17.51041666666667
$ws.Range("P29").Value = 48
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = 1
$ws.Range("S29").Value = 1
$ws.Range("A30").Value = 17
$ws.Range("B30").Value = 'memory'
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 29
$ws.Range("F30").Value = 189
$ws.Range("G30").Value = 'living_rooms'
$ws.Range("H30").Value = 'living_rooms'
$ws.Range("I30").Value = $null
$ws.Range("J30").Value = 'new'
$ws.Range("K30").Value = 'f'
$ws.Range("L30").Value = 'stimuli/img_lgxzn.png'
$ws.Range("M30").Value = 73.11363636363636
$ws.Range("N30").Value = 49.97727272727273
$ws.Range("O30").Value = 61.54545454545455
$ws.Range("P30").Value = 44
$ws.Range("Q30").Value = 6
$ws.Range("R30").Value = 6
$ws.Range("S30").Value = 6
$ws.Range("A31").Value = 17
$ws.Range("B31").Value = 'memory'
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 30
$ws.Range("F31").Value = 190
$ws.Range("G31").Value = 'living_rooms'
$ws.Range("H31").Value = 'living_rooms'
$ws.Range("I31").Value = $null
$ws.Range("J31").Value = 'new'
$ws.Range("K31").Value = 'f'
$ws.Range("L31").Value = 'stimuli/img_pjfx6.png'
$ws.Range("M31").Value = 32.23404255319149
$ws.Range("N31").Value = 26.59574468085106
$ws.Range("O31").Value = 29.41489361702127
$ws.Range("P31").Value = 47
$ws.Range("Q31").Value = 2
$ws.Range("R31").Value = 2
$ws.Range("S31").Value = 2
$ws.Range("A32").Value = 17
$ws.Range("B32").Value = 'memory'
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 31
$ws.Range("F32").Value = 191
$ws.Range("G32").Value = 'living_rooms'
$ws.Range("H32").Value = 'living_rooms'
$ws.Range("I32").Value = $null
$ws.Range("J32").Value = 'new'
$ws.Range("K32").Value = 'f'
$ws.Range("L32").Value = 'stimuli/img_89dvt.png'
$ws.Range("M32").Value = 81.09756097560975
$ws.Range("N32").Value = 64.6829268292683
$ws.Range("O32").Value = 72.89024390243902
$ws.Range("P32").Value = 41
$ws.Range("Q32").Value = 8
$ws.Range("R32").Value = 8
$ws.Range("S32").Value = 8
$ws.Range("A33").Value = 17
$ws.Range("B33").Value = 'memory'
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 32
$ws.Range("F33").Value = 192
$ws.Range("G33").Value = 'living_rooms'
$ws.Range("H33").Value = 'living_rooms'
$ws.Range("I33").Value = 'target'
$ws.Range("J33").Value = 'old'
$ws.Range("K33").Value = 'j'
$ws.Range("L33").Value = 'stimuli/img_abobq.png'
$ws.Range("M33").Value = 75.1842105263158
$ws.Range("N33").Value = 54.13157894736842
$ws.Range("O33").Value = 64.65789473684211
$ws.Range("P33").Value = 38
$ws.Range("Q33").Value = 6
$ws.Range("R33").Value = 6
$ws.Range("S33").Value = 6
$ws.Range("A34").Value = 17
$ws.Range("B34").Value = 'memory'
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 33
$ws.Range("F34").Value = 193
$ws.Range("G34").Value = 'living_rooms'
$ws.Range("H34").Value = 'living_rooms'
$ws.Range("I34").Value = $null
$ws.Range("J34").Value = 'new'
$ws.Range("K34").Value = 'f'
$ws.Range("L34").Value = 'stimuli/img_fmgjx.png'
$ws.Range("M34").Value = 79.9
$ws.Range("N34").Value = 56.975
$ws.Range("O34").Value = 68.4375
$ws.Range("P34").Value = 40
$ws.Range("Q34").Value = 7
$ws.Range("R34").Value = 7
$ws.Range("S34").Value = 7
$ws.Range("A35").Value = 17
$ws.Range("B35").Value = 'memory'
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 34
$ws.Range("F35").Value = 194
$ws.Range("G35").Value = 'living_rooms'
$ws.Range("H35").Value = 'living_rooms'
$ws.Range("I35").Value = 'target'
$ws.Range("J35").Value = 'old'
$ws.Range("K35").Value = 'j'
$ws.Range("L35").Value = 'stimuli/img_xy930.png'
$ws.Range("M35").Value = 70.5952380952381
$ws.Range("N35").Value = 49.47619047619047
$ws.Range("O35").Value = 60.03571428571429
$ws.Range("P35").Value = 42
$ws.Range("Q35").Value = 6
$ws.Range("R35").Value = 6
$ws.Range("S35").Value = 6
$ws.Range("A36").Value = 17
$ws.Range("B36").Value = 'memory'
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 35
$ws.Range("F36").Value = 195
$ws.Range("G36").Value = 'living_rooms'
$ws.Range("H36").Value = 'living_rooms'
$ws.Range("I36").Value = 'target'
$ws.Range("J36").Value = 'old'
$ws.Range("K36").Value = 'j'
$ws.Range("L36").Value = 'stimuli/img_0kqc0.png'
$ws.Range("M36").Value = 43.74468085106383
$ws.Range("N36").Value = 27.14893617021277
$ws.Range("O36").Value = 35.4468085106383
$ws.Range("P36").Value = 47
$ws.Range("Q36").Value = 2
$ws.Range("R36").Value = 2
$ws.Range("S36").Value = 2
$ws.Range("A37").Value = 17
$ws.Range("B37").Value = 'memory'
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 36
$ws.Range("F37").Value = 196
$ws.Range("G37").Value = 'living_rooms'
$ws.Range("H37").Value = 'living_rooms'
$ws.Range("I37").Value = $null
$ws.Range("J37").Value = 'new'
$ws.Range("K37").Value = 'f'
$ws.Range("L37").Value = 'stimuli/img_7lz7m.png'
$ws.Range("M37").Value = 51.5531914893617
$ws.Range("N37").Value = 32.87234042553192
$ws.Range("O37").Value = 42.21276595744681
$ws.Range("P37").Value = 47
$ws.Range("Q37").Value = 3
$ws.Range("R37").Value = 3
$ws.Range("S37").Value = 3
$ws.Range("A38").Value = 17
$ws.Range("B38").Value = 'memory'
$ws.Range("C38").Value = 3
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 37
$ws.Range("F38").Value = 197
$ws.Range("G38").Value = 'living_rooms'
$ws.Range("H38").Value = 'living_rooms'
$ws.Range("I38").Value = $null
$ws.Range("J38").Value = 'new'
$ws.Range("K38").Value = 'f'
$ws.Range("L38").Value = 'stimuli/img_emh91.png'
$ws.Range("M38").Value = 82.06666666666666
$ws.Range("N38").Value = 63.33333333333334
$ws.Range("O38").Value = 72.7
$ws.Range("P38").Value = 45
$ws.Range("Q38").Value = 8
$ws.Range("R38").Value = 8
$ws.Range("S38").Value = 8
$ws.Range("A39").Value = 17
$ws.Range("B39").Value = 'memory'
$ws.Range("C39").Value = 3
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 38
$ws.Range("F39").Value = 198
$ws.Range("G39").Value = 'living_rooms'
$ws.Range("H39").Value = 'living_rooms'
$ws.Range("I39").Value = 'target'
$ws.Range("J39").Value = 'old'
$ws.Range("K39").Value = 'j'
$ws.Range("L39").Value = 'stimuli/img_6a0hu.png'
$ws.Range("M39").Value = 61.275
$ws.Range("N39").Value = 42.025
$ws.Range("O39").Value = 51.65
$ws.Range("P39").Value = 40
$ws.Range("Q39").Value = 4
$ws.Range("R39").Value = 4
$ws.Range("S39").Value = 4
$ws.Range("A40").Value = 17
$ws.Range("B40").Value = 'memory'
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 39
$ws.Range("F40").Value = 199
$ws.Range("G40").Value = 'living_rooms'
$ws.Range("H40").Value = 'living_rooms'
$ws.Range("I40").Value = $null
$ws.Range("J40").Value = 'new'
$ws.Range("K40").Value = 'f'
$ws.Range("L40").Value = 'stimuli/img_lpas9.png'
$ws.Range("M40").Value = 59.36585365853659
$ws.Range("N40").Value = 39.09756097560975
$ws.Range("O40").Value = 49.23170731707317
$ws.Range("P40").Value = 41
$ws.Range("Q40").Value = 4
$ws.Range("R40").Value = 4
$ws.Range("S40").Value = 4
$ws.Range("A41").Value = 17
$ws.Range("B41").Value = 'memory'
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 40
$ws.Range("F41").Value = 200
$ws.Range("G41").Value = 'living_rooms'
$ws.Range("H41").Value = 'living_rooms'
$ws.Range("I41").Value = 'target'
$ws.Range("J41").Value = 'old'
$ws.Range("K41").Value = 'j'
$ws.Range("L41").Value = 'stimuli/img_4o8l0.png'
$ws.Range("M41").Value = 46.02173913043478
$ws.Range("N41").Value = 31.45652173913043
$ws.Range("O41").Value = 38.73913043478261
$ws.Range("P41").Value = 46
$ws.Range("Q41").Value = 3
$ws.Range("R41").Value = 3
$ws.Range("S41").Value = 3
$ws.Range("A42").Value = 17
$ws.Range("B42").Value = 'memory'
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 41
$ws.Range("F42").Value = 201
$ws.Range("G42").Value = 'living_rooms'
$ws.Range("H42").Value = 'living_rooms'
$ws.Range("I42").Value = 'target'
$ws.Range("J42").Value = 'old'
$ws.Range("K42").Value = 'j'
$ws.Range("L42").Value = 'stimuli/img_wz6x5.png'
$ws.Range("M42").Value = 68.3695652173913
$ws.Range("N42").Value = 48.47826086956522
$ws.Range("O42").Value = 58.42391304347826
$ws.Range("P42").Value = 46
$ws.Range("Q42").Value = 5
$ws.Range("R42").Value = 5
$ws.Range("S42").Value = 5
